$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 43; $r -le 1033; $r += 10) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
